# Regenerate save_data: replace column G ("K", strikeouts) values for rows 2-39
# with freshly computed K counts (was previously populated from "Strike#").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(3,5,6,6,5,2,7,8,7,4,8,6,5,2,5,9,6,5,5,7,6,7,2,4,2,3,2,1,1,1,1,0,3,1,2,0,2,2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
